$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last data row (row 1 is header; data goes through last used row)
$lastRow = $ws.UsedRange.Rows.Count

# Add new header cells for the season record columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header style used by the rest of row 1 (bold, centered, bordered)
$headerStyleRange = $ws.Range("AD1:AF1")
$headerStyleRange.Font.Bold = $true
$headerStyleRange.HorizontalAlignment = -4108
$headerStyleRange.VerticalAlignment = -4160
$headerStyleRange.Borders.LineStyle = 1

# Fill in the season record (Wins/Losses/Ties) for every data row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 96
    $ws.Cells.Item($r, 31).Value = 66
    $ws.Cells.Item($r, 32).Value = 0
}

Write-Output "done"
